$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 26 data - copy formatting from row 25 (date style) first
$ws.Range("A25").Copy()
$ws.Range("A26").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A26").Value = 42361
$ws.Range("B26").Value = "80 min"
$ws.Range("C26").Value = "Ronald"
$ws.Range("D26").Value = "V2.0 bijwerken, pagina's toegevoegd"
$ws.Range("E26").Value = "gecontroleerd door John"

# Update selection to reflect where the user ended up (B27)
$ws.Range("B27").Select()
